$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update product names in column B (rows 6-11) with new fruit/vegetable entries.
# B7 and B11 originally carry a quote-prefixed ("Normal" text forced via leading
# apostrophe) style, so re-apply the leading apostrophe to keep that formatting.
$ws.Range("B6").Value = "Grapes_Fruit"
$ws.Range("B7").Value = "'Mango_Fruit"
$ws.Range("B8").Value = "Cabbage_Vegetable"
$ws.Range("B9").Value = "Raddish_Vegetable"
$ws.Range("B10").Value = "Guava_Fruit"
$ws.Range("B11").Value = "'Strawberry_Fruit"

# Clear out the now-unused trailing rows (12-14), keeping their styling
$ws.Range("A12:B14").ClearContents()

# Widen column B to fit the new, longer product names
$ws.Columns.Item(2).ColumnWidth = 16.83

# Move the active selection
$ws.Range("B18").Select()
